$d = $word.ActiveDocument

# Replace the entire body content with the new Chinese bio text, laid out as
# several paragraphs (matching the target structure): four text paragraphs,
# each with a 720-twip first-line indent, separated by a bookmark-only empty
# paragraph and two fully empty paragraphs. The original trailing empty
# paragraph (with the stray Times New Roman run-formatting) is dropped.

$full = $d.Range(0, $d.Content.End)

$xmlSnippet = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:ind w:firstLine="720"/></w:pPr>
<w:r><w:t>Amy C. Edmondson 是哈佛商学院领导力和管理学 Novartis 讲座教授，以及该学院技术与运营管理系副主任。“Novartis 教席”旨在激励人际互动研究，推动建立成功的商业企业，为全人类造福。</w:t></w:r>
</w:p>
<w:p>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
<w:p>
<w:pPr><w:ind w:firstLine="720"/></w:pPr>
<w:r><w:t>Edmondson 讲授领导力、团队决策、组织学习学科的 MBA 课程和高管培训课程以及实地研究方法学科的博士课程。她的研究涉及团队和组织中领导力对学习、协作和创新的影响，曾在学术杂志、管理期刊和书籍中发表了 60 多篇文章。</w:t></w:r>
</w:p>
<w:p/>
<w:p>
<w:pPr><w:ind w:firstLine="720"/></w:pPr>
<w:r><w:t>2003 年，管理学院组织行为系评选 Edmondson 获得“康明奖”，藉此表彰她在职业生涯的早中期取得的突出成就；2000 年评选她的文章《Psychological safety and learning behavior in work teams》获得该领域当年的年度最佳论文奖。她与 Anita Tucker 合著的文章《Why Hospitals Don't Learn from Failures: Organizational and Psychological Dynamics That Inhibit System Change》获得了 2004 年度“埃森哲奖”，藉此表彰她对管理实践的重大贡献。</w:t></w:r>
</w:p>
<w:p/>
<w:p>
<w:pPr><w:ind w:firstLine="720"/></w:pPr>
<w:r><w:t>Edmondson 拥有哈佛大学组织行为学博士学位、心理学硕士学位和工程与设计学学士学位。</w:t></w:r>
</w:p>
<w:sectPr>
<w:pgSz w:w="12240" w:h="15840"/>
<w:pgMar w:top="1440" w:right="1440" w:bottom="1440" w:left="1440" w:header="720" w:footer="720" w:gutter="0"/>
<w:cols w:space="720"/>
<w:docGrid w:linePitch="360"/>
</w:sectPr>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$full.InsertXML($xmlSnippet)
